$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# E49 was stored as text "590024"; change it to a real number.
$ws.Range("E49").Value = 590024

# Append new row 50 with the breakout record.
$ws.Range("A50").Value = "25/06/2024 06:44:42"
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = "FACT"
$ws.Range("D50").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "590024"
$ws.Range("F50").Value = -0.58
$ws.Range("G50").Value = 1015.05
$ws.Range("H50").Value = 1493920
